$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 19) so the table shrinks from 18 to 17 players
$ws.Rows.Item(19).Delete()

# Rewrite the player table (rows 2-18) with the updated roster
$data = @(
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
